$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 197.76923
$ws.Range("I9").Value = 193.2
$ws.Range("K9").Value = 193.2
$ws.Range("M9").Value = -24.19999999999999
$ws.Range("H17").Value = 64789.785
$ws.Range("J17").Value = 64789.785
$ws.Range("L17").Value = 194369.355
$ws.Range("N17").Value = -194705.355
$ws.Range("H43").Value = 6105.3335
$ws.Range("I43").Value = 6458.25
$ws.Range("J43").Value = 5977
$ws.Range("K43").Value = 6458.25
$ws.Range("L43").Value = 5977
$ws.Range("M43").Value = -6389.25
$ws.Range("N43").Value = -6115
$ws.Range("H46").Value = 1430395.9
$ws.Range("I46").Value = 1499.5
$ws.Range("J46").Value = 2001954.4
$ws.Range("K46").Value = 4498.5
$ws.Range("L46").Value = 6005863.199999999
$ws.Range("M46").Value = -4379.5
$ws.Range("N46").Value = -6006101.199999999
$ws.Range("H53").Value = 189.16667
$ws.Range("I53").Value = 203.08333
$ws.Range("J53").Value = 161.33333
$ws.Range("K53").Value = 203.08333
$ws.Range("L53").Value = 161.33333
$ws.Range("M53").Value = 433.91667
$ws.Range("N53").Value = -1435.33333
$ws.Range("H60").Value = 1430395.9
$ws.Range("I60").Value = 1499.5
$ws.Range("J60").Value = 2001954.4
$ws.Range("K60").Value = 4498.5
$ws.Range("L60").Value = 6005863.199999999
$ws.Range("M60").Value = -4014.5
$ws.Range("N60").Value = -6006831.199999999
$ws.Range("H132").Value = 4605.706
$ws.Range("J132").Value = 6300.2856
$ws.Range("L132").Value = 18900.8568
$ws.Range("N132").Value = -23960.8568
$ws.Range("H137").Value = 93087.17999999999
$ws.Range("J137").Value = 251587.25
$ws.Range("L137").Value = 754761.75
$ws.Range("N137").Value = -759861.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8000
$ws.Range("I61").Value = 8000
$ws.Range("K61").Value = 8000
$ws.Range("M61").Value = -7788
$ws.Range("H74").Value = 2106.077
$ws.Range("I74").Value = 2134.4546
$ws.Range("K74").Value = 2134.4546
$ws.Range("M74").Value = -1260.4546
$ws.Range("H77").Value = 2106.077
$ws.Range("I77").Value = 2134.4546
$ws.Range("K77").Value = 10672.273
$ws.Range("M77").Value = -6304.273000000001
$ws.Range("H107").Value = 30000
$ws.Range("J107").Value = 30000
$ws.Range("L107").Value = 30000
$ws.Range("N107").Value = -37680
$ws.Range("H122").Value = 20476.273
$ws.Range("I122").Value = 30463
$ws.Range("K122").Value = 91389
$ws.Range("M122").Value = -88939
$ws.Range("H136").Value = 8000
$ws.Range("I136").Value = 8000
$ws.Range("K136").Value = 24000
$ws.Range("M136").Value = -21450

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2312.0227
$ws.Range("I31").Value = 1839.2593
$ws.Range("J31").Value = 3062.8823
$ws.Range("K31").Value = 1839.2593
$ws.Range("L31").Value = 3062.8823
$ws.Range("M31").Value = -1544.2593
$ws.Range("N31").Value = -3652.8823
$ws.Range("H34").Value = 2312.0227
$ws.Range("I34").Value = 1839.2593
$ws.Range("J34").Value = 3062.8823
$ws.Range("K34").Value = 1839.2593
$ws.Range("L34").Value = 3062.8823
$ws.Range("M34").Value = -1637.2593
$ws.Range("N34").Value = -3466.8823
$ws.Range("H59").Value = 14157.2
$ws.Range("I59").Value = 9552
$ws.Range("K59").Value = 9552
$ws.Range("M59").Value = -8407
$ws.Range("H107").Value = 1888.7778
$ws.Range("I107").Value = 1648
$ws.Range("K107").Value = 1648
$ws.Range("M107").Value = 272

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 32.9375
$ws.Range("I10").Value = 28.466667
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 85.400001
$ws.Range("L10").Value = 300
$ws.Range("M10").Value = 53.599999
$ws.Range("N10").Value = -578
$ws.Range("H11").Value = 36121.25
$ws.Range("I11").Value = 424.82608
$ws.Range("K11").Value = 1274.47824
$ws.Range("M11").Value = -1134.47824
$ws.Range("H12").Value = 799.55554
$ws.Range("J12").Value = 274.5
$ws.Range("L12").Value = 823.5
$ws.Range("N12").Value = -1169.5
$ws.Range("H13").Value = 110
$ws.Range("I13").Value = 76
$ws.Range("J13").Value = 166.66667
$ws.Range("K13").Value = 228
$ws.Range("L13").Value = 500.00001
$ws.Range("M13").Value = -60
$ws.Range("N13").Value = -836.00001
$ws.Range("H26").Value = 2261
$ws.Range("I26").Value = 153.8
$ws.Range("J26").Value = 4368.2
$ws.Range("K26").Value = 461.4
$ws.Range("L26").Value = 13104.6
$ws.Range("M26").Value = -173.4
$ws.Range("N26").Value = -13680.6
$ws.Range("H51").Value = 1473.25
$ws.Range("I51").Value = 1473.25
$ws.Range("J51").Value = 0
$ws.Range("K51").Value = 4419.75
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = -3959.75
$ws.Range("N51").ClearContents()
$ws.Range("H57").Value = 990
$ws.Range("I57").Value = 990
$ws.Range("K57").Value = 2970
$ws.Range("M57").Value = -2411
$ws.Range("H60").Value = 841.4167
$ws.Range("I60").Value = 442.42856
$ws.Range("J60").Value = 1400
$ws.Range("K60").Value = 1327.28568
$ws.Range("L60").Value = 4200
$ws.Range("M60").Value = -1076.28568
$ws.Range("N60").Value = -4702
$ws.Range("H81").Value = 35722052
$ws.Range("I81").Value = 3299
$ws.Range("J81").Value = 62511116
$ws.Range("K81").Value = 9897
$ws.Range("L81").Value = 187533348
$ws.Range("M81").Value = -8774
$ws.Range("N81").Value = -187535594
$ws.Range("H84").Value = 35722052
$ws.Range("I84").Value = 3299
$ws.Range("J84").Value = 62511116
$ws.Range("K84").Value = 29691
$ws.Range("L84").Value = 562600044
$ws.Range("M84").Value = -24075
$ws.Range("N84").Value = -562611276
$ws.Range("H113").Value = 478.4
$ws.Range("I113").Value = 348.125
$ws.Range("J113").Value = 999.5
$ws.Range("K113").Value = 1044.375
$ws.Range("L113").Value = 2998.5
$ws.Range("M113").Value = 1125.625
$ws.Range("N113").Value = -7338.5
$ws.Range("H115").Value = 1599.8334
$ws.Range("I115").Value = 899.75
$ws.Range("K115").Value = 2699.25
$ws.Range("M115").Value = -1524.25
$ws.Range("H117").Value = 5047.1875
$ws.Range("J117").Value = 6610.5835
$ws.Range("L117").Value = 19831.7505
$ws.Range("N117").Value = -26715.7505
$ws.Range("H132").Value = 2073.125
$ws.Range("I132").Value = 1235.2
$ws.Range("J132").Value = 2454
$ws.Range("K132").Value = 11116.8
$ws.Range("L132").Value = 22086
$ws.Range("M132").Value = -8586.800000000001
$ws.Range("N132").Value = -27146
$ws.Range("H140").Value = 3152.6667
$ws.Range("J140").Value = 10833.25
$ws.Range("L140").Value = 32499.75
$ws.Range("N140").Value = -42859.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2714.2
$ws.Range("I102").Value = 2434.7144
$ws.Range("J102").Value = 3366.3333
$ws.Range("K102").Value = 2434.7144
$ws.Range("L102").Value = 3366.3333
$ws.Range("M102").Value = -812.7143999999998
$ws.Range("N102").Value = -6610.3333
$ws.Range("H109").Value = 100285
$ws.Range("J109").Value = 100285
$ws.Range("L109").Value = 100285
$ws.Range("N109").Value = -102365
$ws.Range("H126").Value = 3294
$ws.Range("I126").Value = 3349.2856
$ws.Range("J126").Value = 3197.25
$ws.Range("K126").Value = 10047.8568
$ws.Range("L126").Value = 9591.75
$ws.Range("M126").Value = -7577.856800000001
$ws.Range("N126").Value = -14531.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1556.2727
$ws.Range("J16").Value = 400
$ws.Range("L16").Value = 400
$ws.Range("N16").Value = -740
$ws.Range("H100").Value = 24727
$ws.Range("I100").Value = 4201.25
$ws.Range("K100").Value = 4201.25
$ws.Range("M100").Value = -3660.25
$ws.Range("H132").Value = 9792.733
$ws.Range("I132").Value = 13566.777
$ws.Range("J132").Value = 4131.6665
$ws.Range("K132").Value = 40700.331
$ws.Range("L132").Value = 12394.9995
$ws.Range("M132").Value = -38170.331
$ws.Range("N132").Value = -17454.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 559.08
$ws.Range("I113").Value = 304.16666
$ws.Range("K113").Value = 912.4999799999999
$ws.Range("M113").Value = 1257.50002
$ws.Range("H126").Value = 2380.8235
$ws.Range("I126").Value = 1772.091
$ws.Range("K126").Value = 5316.272999999999
$ws.Range("M126").Value = -2846.272999999999
$ws.Range("H132").Value = 40572.285
$ws.Range("I132").Value = 40572.285
$ws.Range("K132").Value = 121716.855
$ws.Range("M132").Value = -119186.855
$ws.Range("H136").Value = 5759.077
$ws.Range("I136").Value = 6156.8
$ws.Range("J136").Value = 4433.3335
$ws.Range("K136").Value = 18470.4
$ws.Range("L136").Value = 13300.0005
$ws.Range("M136").Value = -15920.4
$ws.Range("N136").Value = -18400.0005
